# Upload Karyawan template - tambah kolom cuti (Jatah Cuti Pribadi, Jatah Cuti
# Bersama, Jatah Cuti Tahun Lalu, Expired Date Cuti Tahun Lalu, Hutang Cuti)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Salin format kolom AD (header + contoh) ke kolom-kolom baru AE:AI supaya
# style-nya (font/fill/border) sama persis dengan header/contoh yang sudah ada.
$ws.Range("AD1:AD2").Copy()
$ws.Range("AE1:AI2").PasteSpecial(-4122)

# Header baris 1
$ws.Range("AE1").Value = "Jatah Cuti Pribadi *"
$ws.Range("AF1").Value = "Jatah Cuti Bersama *"
$ws.Range("AG1").Value = "Jatah Cuti Tahun Lalu *"

# Contoh isian baris 2
$ws.Range("AE2").Value = "Isi dengan angka saja"
$ws.Range("AF2").Value = "Isi dengan angka saja"
$ws.Range("AG2").Value = "Isi dengan angka saja"
$ws.Range("AH2").Value = "Ex : 01/02/2025"

$ws.Range("AH1").Value = "Expired Date Cuti Tahun Lalu *"
$ws.Range("AI1").Value = "Hutang Cuti *"
$ws.Range("AI2").Value = "Isi dengan angka saja"

# Lebar kolom baru mengikuti lebar konten (mirip AutoFit manual)
$ws.Columns("AE").ColumnWidth = 23.140625
$ws.Columns("AF").ColumnWidth = 24.85546875
$ws.Columns("AG").ColumnWidth = 27.140625
$ws.Columns("AH").ColumnWidth = 36.7109375
$ws.Columns("AI").ColumnWidth = 29.7109375

# Golongan Darah (Z2) tidak perlu highlight style lagi - samakan dengan sel
# contoh lain yang polos (tanpa style index khusus)
$ws.Range("A2").Copy()
$ws.Range("Z2").PasteSpecial(-4122)

# Pindahkan sorotan/selection ke area kolom baru
$ws.Range("AI4").Select()

Write-Host "Upload karyawan template updated with Cuti columns (AE:AI)."
